$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 709 — shifts the existing rows 709:750
# (and their values) down to 710:751, extending the used range to D751.
$ws.Rows.Item(709).Insert()

# Populate the newly-inserted row 709 with the new daily entry
# (2026/01/25, Sunday, hour 8, ranking 201). Force text formatting before
# assigning the date-like string so it is stored as literal text (matching
# every other date cell in column A) instead of being auto-parsed into an
# Excel date serial number; ClearFormats afterwards drops the now-unneeded
# "@" number format so the cell ends up with the default (unstyled) look,
# same as its neighbours.
$cellA = $ws.Cells.Item(709, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026/01/25"
$cellA.ClearFormats()

$ws.Cells.Item(709, 2).Value = "日"
$ws.Cells.Item(709, 3).Value = 8
$ws.Cells.Item(709, 4).Value = 201
